$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row height corrections (rows shrank by one wrapped line each) ---
$ws.Rows.Item(3).RowHeight = 15.75
$ws.Rows.Item(4).RowHeight = 15.75
$ws.Rows.Item(5).RowHeight = 15.75
$ws.Rows.Item(6).RowHeight = 15.75
$ws.Rows.Item(7).RowHeight = 30.75
$ws.Rows.Item(8).RowHeight = 30.75
$ws.Rows.Item(9).RowHeight = 15.75
$ws.Rows.Item(13).RowHeight = 30.75
$ws.Rows.Item(14).RowHeight = 15.75
$ws.Rows.Item(16).RowHeight = 30.75
$ws.Rows.Item(18).RowHeight = 15.75

# --- Row 23: fill in missing Heart-Rate exercise readings & correct values ---
$ws.Cells.Item(23,2).Value2 = 79
$ws.Cells.Item(23,3).Value2 = 124
$ws.Cells.Item(23,4).Value2 = 136
$ws.Cells.Item(23,5).Value2 = 145
$ws.Cells.Item(23,6).Value2 = 157
$ws.Cells.Item(23,7).Value2 = 174
$ws.Range("B23:G23").HorizontalAlignment = 1

# --- Row 25 / 26: corrected totals ---
$ws.Cells.Item(25,2).Value2 = 0.20833333333333334
$ws.Cells.Item(26,2).Value2 = 2640

# --- Insert a note row above the "Speed" exercise table ---
$ws.Rows.Item(27).Insert()
$ws.Rows.Item(27).Insert()
$ws.Rows.Item(27).Clear()
$ws.Rows.Item(28).Clear()

$ws.Cells.Item(27,1).Value2 = "***No reset used between the first chart and exercise and heart rate at 0 taken with patient standing. Patient didn’t give up"
$ws.Range("A27:L27").Interior.Color = 65535
$ws.Cells.Item(27,1).HorizontalAlignment = -4131
$ws.Cells.Item(27,1).IndentLevel = 5
$ws.Rows.Item(27).RowHeight = 15.75

# --- Update view: scroll position and selection ---
$ws.Range("A27:L27").Select()

Write-Host "done"
